$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Remove the empty/unused "title" placeholder shape (id 206,
# "Google Shape;206;p7") that was left over with no text content.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Google Shape;206;p7") {
        $sh.Delete()
    }
}
